$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(3, 4).Value = 'Background
id="Par1">COVID-19 is a virus pandemic.

 According to the first obtained data, COVID-19 has defined with findings such as cough, fever, diarrhea, and fatigue although neurological symptoms of patients with COVID-19 have not been investigated in detail.

 This study aims to investigate the neurological findings via obtained face-to-face anamnesis and detailed neurological examination in patients with COVID-19.
Methods
id="Par2">Two hundred thirty-nine consecutive inpatients with COVID-19, supported with laboratory tests, were evaluated.

 Detailed neurological examinations and evaluations of all patients were performed.

 All evaluations and examinations were performed by two neurologists who have at least five-year experience.


Results
id="Par3">This study was carried out 239 patients (133 male + 106 female) with diagnosed COVID-19. Neurological findings were present in 83 of 239 patients (34.7%).

 The most common neurological finding was a headache (27.6%).

 D-dimer blood levels were detected to be significantly higher in patients with at least one neurological symptom than patients without the neurological symptom (p &lt; 0.05).

 IL-6 level was found to be significantly higher in patients with headache than without headache (p &lt; 0.05).

 Creatine kinase (CK) level was detected to be significantly higher in patients with muscle pain (p &lt; 0.05).


Conclusion
id="Par4">Neurological symptoms are often seen in patients with COVID-19. Headache was the most common seen neurological symptom in this disease.

 Dizziness, impaired consciousness, smell and gustation impairments, cerebrovascular disorders, epileptic seizures, and myalgia were detected as other findings apart from the headache.

 It is suggested that determining these neurological symptoms prevents the diagnosis delay and helps to prohibit virus spread.


'
$ws.Cells.Item(3, 5).Value = '[Ömer%Karadaş%NULL%1,    Bilgin%Öztürk%Drbilgin@gmail.com%2,    Bilgin%Öztürk%Drbilgin@gmail.com%0,    Ali Rıza%Sonkaya%NULL%2,    Ali Rıza%Sonkaya%NULL%0]'

$ws.Rows(3).AutoFit()
